$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume figures per latest data refresh
$ws.Range("D2").Value = "25.979.66"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "1.593.93"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'210.63"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("D8").Value = "'0.246"
$ws.Range("E8").Value = "  -0.83%  "
$ws.Range("E9").Value = "  -1.30%  "
$ws.Range("E10").Value = "  -1.83%  "
$ws.Range("D11").Value = "'0.0810"
$ws.Range("E11").Value = "  +2.71%  "
$ws.Range("D12").Value = "1.815.88"
$ws.Range("E12").Value = "  +0.27%  "
$ws.Range("D13").Value = "1.596.37"
$ws.Range("E13").Value = "  +0.54%  "
$ws.Range("E14").Value = "  -0.91%  "
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("D16").Value = "25.988.70"
$ws.Range("E16").Value = "  +0.33%  "
$ws.Range("D17").Value = "'60.04"
$ws.Range("E17").Value = "  -0.38%  "
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").Value = "'199.95"
$ws.Range("E20").Value = "  +3.27%  "
$ws.Range("E21").Value = "  +0.55%  "
$ws.Range("D22").Value = "'9.24"
$ws.Range("E22").Value = "  -1.91%  "
$ws.Range("D23").Value = "'5.99"
$ws.Range("E23").Value = "  +0.89%  "
$ws.Range("E24").Value = "  +7.29%  "
$ws.Range("D25").Value = "'143.19"
$ws.Range("E25").Value = "  +1.07%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "'0.121"
$ws.Range("E27").Value = "  -8.47%  "
$ws.Range("E29").Value = "  -0.63%  "
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("D32").Value = "'3.13"
$ws.Range("E32").Value = "  +0.09%  "
$ws.Range("E33").Value = "  -3.36%  "
$ws.Range("E34").Value = "  -1.50%  "
$ws.Range("E35").Value = "  +0.26%  "
$ws.Range("D36").Value = "1.122.92"
$ws.Range("E36").Value = "  +1.31%  "
$ws.Range("D37").Value = "'0.0162"
$ws.Range("E37").Value = "  +7.79%  "
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("E39").Value = "  -1.10%  "
$ws.Range("E40").Value = "  -0.32%  "
$ws.Range("D41").Value = "'0.491"
$ws.Range("E42").Value = "  -4.19%  "
$ws.Range("D43").Value = "1.727.43"
$ws.Range("E43").Value = "  +0.20%  "
$ws.Range("D44").Value = "'92.55"
$ws.Range("E44").Value = "  -1.25%  "
$ws.Range("E45").Value = "  -1.39%  "
$ws.Range("E46").Value = "  -1.39%  "
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("E48").Value = "  -1.44%  "
$ws.Range("D49").Value = "'0.407"
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.15"
$ws.Range("E51").Value = "  -2.32%  "

Write-Output "Updated cryptos list"
